$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (existing rows 9-19 shift down to 10-20)
$ws.Rows(9).Insert()

# The freshly inserted row gets a brand-new (blank) style; restore the
# correct look-and-feel by copying the formatting used by the other rows
# in the first "group" of the table (rows 2-8, style ids 3/10/6).
$ws.Range("A2:J2").Copy()
$ws.Range("A9:J9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the data for the new project row
$ws.Range("A9").Value = "elk"
$ws.Range("B9").Value = "Apontamentos e exemplos relacionados a plataforma Elastic Stack"
$ws.Range("C9").Value = (Get-Date -Year 2017 -Month 7 -Day 7).Date
$ws.Range("D9").Value = "DEV"
$ws.Range("E9").Value = "X"
$ws.Range("H9").Value = "MIT"
$ws.Range("I9").Value = "MPS"

# Update the view so that it matches the saved selection/scroll position
$ws.Range("A3").Select()

"Done"
